# Decrement the "剩余" (remaining) values in column E for all data rows
# (rows 2-99), except row 36 which keeps its original value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 99; $row++) {
    if ($row -eq 36) {
        continue
    }
    $cell = $ws.Cells.Item($row, 5)  # Column E is the 5th column
    $current = $cell.Value2
    if ($current -ne $null) {
        $cell.Value = $current - 1
    }
}
